$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1363.7
$ws.Range("J2").Value = 1107
$ws.Range("L2").Value = 1107
$ws.Range("N2").Value = -1333
$ws.Range("H40").Value = 4032.6667
$ws.Range("I40").Value = 1200
$ws.Range("J40").Value = 6865.3335
$ws.Range("K40").Value = 1200
$ws.Range("L40").Value = 6865.3335
$ws.Range("M40").Value = -1025
$ws.Range("N40").Value = -7215.3335
$ws.Range("H112").Value = 4996.619
$ws.Range("I112").Value = 415
$ws.Range("J112").Value = 5478.8945
$ws.Range("K112").Value = 1245
$ws.Range("L112").Value = 16436.6835
$ws.Range("M112").Value = -137
$ws.Range("N112").Value = -18652.6835
$ws.Range("H116").Value = 8508.684999999999
$ws.Range("J116").Value = 3782.8572
$ws.Range("L116").Value = 3782.8572
$ws.Range("N116").Value = -10666.8572
$ws.Range("H131").Value = 2152.8125
$ws.Range("I131").Value = 850
$ws.Range("J131").Value = 2453.4614
$ws.Range("K131").Value = 2550
$ws.Range("L131").Value = 7360.3842
$ws.Range("M131").Value = 2490
$ws.Range("N131").Value = -17440.3842
$ws.Range("H135").Value = 17858538
$ws.Range("I135").Value = 22728228
$ws.Range("J135").Value = 3007
$ws.Range("K135").Value = 204554052
$ws.Range("L135").Value = 27063
$ws.Range("M135").Value = -204551517
$ws.Range("N135").Value = -32133
$ws.Range("H141").Value = 3459.7144
$ws.Range("I141").Value = 1531.6111
$ws.Range("J141").Value = 15028.333
$ws.Range("K141").Value = 4594.8333
$ws.Range("L141").Value = 45084.999
$ws.Range("M141").Value = 585.1666999999998
$ws.Range("N141").Value = -55444.999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5630062
$ws.Range("I32").Value = 6504029.5
$ws.Range("J32").Value = 22104.166
$ws.Range("K32").Value = 6504029.5
$ws.Range("L32").Value = 22104.166
$ws.Range("M32").Value = -6503742.5
$ws.Range("N32").Value = -22678.166
$ws.Range("H61").Value = 22226060
$ws.Range("I61").Value = 27780242
$ws.Range("J61").Value = 9333.333000000001
$ws.Range("K61").Value = 27780242
$ws.Range("L61").Value = 9333.333000000001
$ws.Range("M61").Value = -27780030
$ws.Range("N61").Value = -9757.333000000001
$ws.Range("H97").Value = 760.7143
$ws.Range("I97").Value = 760.7143
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 760.7143
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -264.7143
$ws.Range("N97").ClearContents()
$ws.Range("H104").Value = 77306.25
$ws.Range("J104").Value = 77306.25
$ws.Range("L104").Value = 77306.25
$ws.Range("N104").Value = -84294.25
$ws.Range("H124").Value = 89775.57000000001
$ws.Range("J124").Value = 89775.57000000001
$ws.Range("L124").Value = 89775.57000000001
$ws.Range("N124").Value = -99595.57000000001
$ws.Range("H125").Value = 98715
$ws.Range("J125").Value = 98715
$ws.Range("L125").Value = 98715
$ws.Range("N125").Value = -108555
$ws.Range("H132").Value = 1306405.9
$ws.Range("I132").Value = 1903.5491
$ws.Range("J132").Value = 9622609
$ws.Range("K132").Value = 5710.6473
$ws.Range("L132").Value = 28867827
$ws.Range("M132").Value = -3180.6473
$ws.Range("N132").Value = -28872887
$ws.Range("H136").Value = 22226060
$ws.Range("I136").Value = 27780242
$ws.Range("J136").Value = 9333.333000000001
$ws.Range("K136").Value = 83340726
$ws.Range("L136").Value = 27999.999
$ws.Range("M136").Value = -83338176
$ws.Range("N136").Value = -33099.999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 80000
$ws.Range("J42").Value = 80000
$ws.Range("L42").Value = 80000
$ws.Range("N42").Value = -80656
$ws.Range("H94").Value = 374.83334
$ws.Range("I94").Value = 374.83334
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 374.83334
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = 76.16665999999998
$ws.Range("N94").ClearContents()
$ws.Range("H134").Value = 3681.8635
$ws.Range("I134").Value = 3631.6316
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 10894.8948
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -8359.8948
$ws.Range("N134").Value = -17070
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5693.143
$ws.Range("I31").Value = 1705.625
$ws.Range("J31").Value = 7773.587
$ws.Range("K31").Value = 1705.625
$ws.Range("L31").Value = 7773.587
$ws.Range("M31").Value = -1410.625
$ws.Range("N31").Value = -8363.587
$ws.Range("H34").Value = 5693.143
$ws.Range("I34").Value = 1705.625
$ws.Range("J34").Value = 7773.587
$ws.Range("K34").Value = 1705.625
$ws.Range("L34").Value = 7773.587
$ws.Range("M34").Value = -1503.625
$ws.Range("N34").Value = -8177.587
$ws.Range("H132").Value = 19049778
$ws.Range("I132").Value = 21741010
$ws.Range("J132").Value = 13891583
$ws.Range("K132").Value = 65223030
$ws.Range("L132").Value = 41674749
$ws.Range("M132").Value = -65220500
$ws.Range("N132").Value = -41679809
$ws.Range("H134").Value = 10006557
$ws.Range("I134").Value = 11369906
$ws.Range("J134").Value = 8666.666999999999
$ws.Range("K134").Value = 34109718
$ws.Range("L134").Value = 26000.001
$ws.Range("M134").Value = -34107183
$ws.Range("N134").Value = -31070.001
$ws.Range("H141").Value = 70047.7
$ws.Range("J141").Value = 66437.914
$ws.Range("L141").Value = 66437.914
$ws.Range("N141").Value = -76797.914
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 8561.200000000001
$ws.Range("I81").Value = 1778.25
$ws.Range("J81").Value = 13083.167
$ws.Range("K81").Value = 5334.75
$ws.Range("L81").Value = 39249.501
$ws.Range("M81").Value = -4211.75
$ws.Range("N81").Value = -41495.501
$ws.Range("H84").Value = 8561.200000000001
$ws.Range("I84").Value = 1778.25
$ws.Range("J84").Value = 13083.167
$ws.Range("K84").Value = 16004.25
$ws.Range("L84").Value = 117748.503
$ws.Range("M84").Value = -10388.25
$ws.Range("N84").Value = -128980.503
$ws.Range("H137").Value = 34318.145
$ws.Range("I137").Value = 7363.8945
$ws.Range("J137").Value = 66326.31
$ws.Range("K137").Value = 22091.6835
$ws.Range("L137").Value = 198978.93
$ws.Range("M137").Value = -16991.6835
$ws.Range("N137").Value = -209178.93
$ws.Range("H140").Value = 1961.7241
$ws.Range("I140").Value = 1343.3334
$ws.Range("J140").Value = 2624.2856
$ws.Range("K140").Value = 4030.0002
$ws.Range("L140").Value = 7872.8568
$ws.Range("M140").Value = 1149.9998
$ws.Range("N140").Value = -18232.8568
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5562.2856
$ws.Range("I70").Value = 5511.5557
$ws.Range("J70").Value = 5866.6665
$ws.Range("K70").Value = 5511.5557
$ws.Range("L70").Value = 5866.6665
$ws.Range("M70").Value = -5241.5557
$ws.Range("N70").Value = -6406.6665
$ws.Range("H73").Value = 5562.2856
$ws.Range("I73").Value = 5511.5557
$ws.Range("J73").Value = 5866.6665
$ws.Range("K73").Value = 5511.5557
$ws.Range("L73").Value = 5866.6665
$ws.Range("M73").Value = -4575.5557
$ws.Range("N73").Value = -7738.6665
$ws.Range("H122").Value = 1935.6666
$ws.Range("I122").Value = 1588
$ws.Range("K122").Value = 4764
$ws.Range("M122").Value = -2314
$ws.Range("H132").Value = 50008230
$ws.Range("I132").Value = 83344140
$ws.Range("J132").Value = 4365.625
$ws.Range("K132").Value = 250032420
$ws.Range("L132").Value = 13096.875
$ws.Range("M132").Value = -250029890
$ws.Range("N132").Value = -18156.875
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 16658.285
$ws.Range("J93").Value = 2768
$ws.Range("L93").Value = 2768
$ws.Range("N93").Value = -5264
$ws.Range("H108").Value = 67084
$ws.Range("J108").Value = 67084
$ws.Range("L108").Value = 67084
$ws.Range("N108").Value = -74764
$ws.Range("H127").Value = 30000.033
$ws.Range("J127").Value = 30000.033
$ws.Range("L127").Value = 30000.033
$ws.Range("N127").Value = -39920.033
$ws.Range("H132").Value = 2541.7666
$ws.Range("I132").Value = 2320.625
$ws.Range("J132").Value = 3426.3333
$ws.Range("K132").Value = 6961.875
$ws.Range("L132").Value = 10278.9999
$ws.Range("M132").Value = -4431.875
$ws.Range("N132").Value = -15338.9999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 9944
$ws.Range("J54").Value = 9944
$ws.Range("L54").Value = 9944
$ws.Range("N54").Value = -10984
$ws.Range("H81").Value = 3274.4707
$ws.Range("I81").Value = 3605.6667
$ws.Range("J81").Value = 2479.6
$ws.Range("K81").Value = 7211.3334
$ws.Range("L81").Value = 4959.2
$ws.Range("M81").Value = -6150.3334
$ws.Range("N81").Value = -7081.2
$ws.Range("H84").Value = 3274.4707
$ws.Range("I84").Value = 3605.6667
$ws.Range("J84").Value = 2479.6
$ws.Range("K84").Value = 36056.667
$ws.Range("L84").Value = 24796
$ws.Range("M84").Value = -30752.667
$ws.Range("N84").Value = -35404
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H137").Value = 49933.332
$ws.Range("J137").Value = 49933.332
$ws.Range("L137").Value = 49933.332
$ws.Range("N137").Value = -60133.332
